# "Generate Report for Archive"
#
# The localization run for rows 5-7 (the .yml/.md files handled in the
# 22f3c79d... batch) has moved on from "Ready for handoff" to
# "In Translation" in every place that status shows up:
#   - Overview sheet: per-language status columns (zh-cn / de-de)
#   - zh-cn sheet (table): Status column
#   - de-de sheet (table): Status column
#
# Once nothing references the "Ready for handoff" string any more, Excel
# drops it from the shared-string table on save, which is what shrinks the
# status columns (their autofit width was sized for that longer string).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
foreach ($r in 5..7) {
    $ov.Range("E$r").Value = "In Translation"
    $ov.Range("F$r").Value = "In Translation"
}

# ---- zh-cn sheet (table "zh_cn") -------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
foreach ($r in 5..7) {
    $zh.Range("C$r").Value = "In Translation"
}

# ---- de-de sheet (table "de_de") -------------------------------------
$de = $wb.Worksheets.Item("de-de")
foreach ($r in 5..7) {
    $de.Range("C$r").Value = "In Translation"
}

# ---- Narrow the status columns to match the new (shorter) content ----
# Raw OOXML width 17.2159881591797 -> 13.4101845877511 is ColumnWidth
# 16.33 -> 12.576851254417766 (width = ColumnWidth + 5/6 in this engine).
$newWidth = 12.576851254417766
$ov.Columns.Item(5).ColumnWidth = $newWidth
$ov.Columns.Item(6).ColumnWidth = $newWidth
$zh.Columns.Item(3).ColumnWidth = $newWidth
$de.Columns.Item(3).ColumnWidth = $newWidth
